# Swap the contents of columns D (group-name) and E (group-code)
# for every row in the used range. This reproduces the effect of the
# upstream codeforIATI codelist change where the "group-code" and
# "group-name" columns were swapped in the source data (columns D/E),
# while the underlying "name" column (B) keeps the same displayed value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2

    $dCell.Value2 = $eVal
    $eCell.Value2 = $dVal
}
